$wb = $excel.ActiveWorkbook

# --- Copy the E:H block from rows 23-33 into rows 2-12 on Sheet1 ---
$ws1 = $wb.Worksheets.Item("Sheet1")
$src = $ws1.Range("E23:H33")
$dst = $ws1.Range("E2:H12")
$src.Copy($dst)

# --- Reorder worksheet tabs: Sheet1, Sheet3, Sheet4, Sheet2 ---
$ws1.Move($wb.Worksheets.Item(1))
$ws3 = $wb.Worksheets.Item("Sheet3")
$ws3.Move($null, $ws1)
$ws4 = $wb.Worksheets.Item("Sheet4")
$ws4.Move($null, $ws3)

# --- Make Sheet1 the active sheet and set the view/selection ---
$ws1.Activate()
$ws1.Range("A5").Select()
$excel.ActiveWindow.ScrollRow = 5
$ws1.Range("E13").Select()
